# Update the income-statement data to the latest reporting period:
# drop the oldest "1396/12" column, shift D:H one column to the left,
# and populate the new "1401/12" column (H) with the freshly published figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial-period column headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"

# --- Row 9: publish-date column headers ---
$ws.Range("D9").Value = "1399-02-31 (8)"
$ws.Range("E9").Value = "1400-04-05 (11)"
$ws.Range("F9").Value = "1401-04-16 (8)"

$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("G9").Value = "1402-02-28 (8)"

# H9 ("1402-02-28") looks like an ISO date, so a plain assignment gets
# auto-converted to a date serial by Excel's input parser. Force it to
# stay literal text, then restore the original (General/no-format) cell
# style by pulling formatting back from the identically-styled G9 cell.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-28"
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 11-27: shift the five yearly-figure columns left by one and fill in the new "1401/12" column (H) ---
# Row 11
$ws.Range("D11").Value = 3258510
$ws.Range("E11").Value = 9728600
$ws.Range("F11").Value = 18982765
$ws.Range("G11").Value = 22970469
$ws.Range("H11").Value = 46336362

# Row 12
$ws.Range("D12").Value = -2659552
$ws.Range("E12").Value = -5672944
$ws.Range("F12").Value = -10519343
$ws.Range("G12").Value = -15273147
$ws.Range("H12").Value = -35270207

# Row 13
$ws.Range("D13").Value = 598958
$ws.Range("E13").Value = 4055656
$ws.Range("F13").Value = 8463422
$ws.Range("G13").Value = 7697322
$ws.Range("H13").Value = 11066155

# Row 14
$ws.Range("D14").Value = -342743
$ws.Range("E14").Value = -1145291
$ws.Range("F14").Value = -1695215
$ws.Range("G14").Value = -2366989
$ws.Range("H14").Value = -3256229

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# Row 16
$ws.Range("D16").Value = 65069
$ws.Range("E16").Value = 105010
$ws.Range("F16").Value = 206804
$ws.Range("G16").Value = 150192
$ws.Range("H16").Value = 1063398

# Row 17
$ws.Range("D17").Value = 321284
$ws.Range("E17").Value = 3015375
$ws.Range("F17").Value = 6975011
$ws.Range("G17").Value = 5480525
$ws.Range("H17").Value = 8873324

# Row 18
$ws.Range("D18").Value = -249063
$ws.Range("E18").Value = -573607
$ws.Range("F18").Value = -480711
$ws.Range("G18").Value = -528202
$ws.Range("H18").Value = -797485

# Row 19
$ws.Range("D19").Value = 82682
$ws.Range("E19").Value = 299162
$ws.Range("F19").Value = 67460
$ws.Range("G19").Value = 100484
$ws.Range("H19").Value = 103002

# Row 20
$ws.Range("D20").Value = 154903
$ws.Range("E20").Value = 2740930
$ws.Range("F20").Value = 6561760
$ws.Range("G20").Value = 5052807
$ws.Range("H20").Value = 8178841

# Row 21
$ws.Range("D21").Value = -60002
$ws.Range("E21").Value = -563867
$ws.Range("F21").Value = -1029317
$ws.Range("G21").Value = -68848
$ws.Range("H21").Value = -69821

# Row 22
$ws.Range("D22").Value = 94901
$ws.Range("E22").Value = 2177063
$ws.Range("F22").Value = 5532443
$ws.Range("G22").Value = 4983959
$ws.Range("H22").Value = 8109020

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# Row 24
$ws.Range("D24").Value = 94901
$ws.Range("E24").Value = 2177063
$ws.Range("F24").Value = 5532443
$ws.Range("G24").Value = 4983959
$ws.Range("H24").Value = 8109020

# Row 25
$ws.Range("D25").Value = 73
$ws.Range("E25").Value = 1361
$ws.Range("F25").Value = 692
$ws.Range("G25").Value = 277
$ws.Range("H25").Value = 451

# Row 26
$ws.Range("D26").Value = 1300000
$ws.Range("E26").Value = 1600000
$ws.Range("F26").Value = 8000000
$ws.Range("G26").Value = 18000000
$ws.Range("H26").Value = 18000000

# Row 27
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 121
$ws.Range("F27").Value = 307
$ws.Range("G27").Value = 277
$ws.Range("H27").Value = 451
